# "data till 8Feb 8PM" - enter the day's collection figures into column N
# (8-Feb-2021) for every customer who paid that day. The F-column
# per-row totals (SUM(G:AK)) and the row-2 summary totals are formulas,
# so they recalculate automatically once the new inputs are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> amount collected on 8-Feb (column N) for this update.
# Rows marked in $highlighted additionally get the purple "entered"
# fill that the workbook already uses to flag certain paid cells.
$amounts = [ordered]@{
    5  = 1000
    6  = 3000
    7  = 8000
    13 = 5000
    14 = 4000
    15 = 2000
    17 = 1000
    22 = 1000
    23 = 3000
    24 = 5000
    29 = 5000
    31 = 1500
    32 = 1000
    33 = 1000
    34 = 2500
    39 = 4000
    41 = 3000
    44 = 3000
    46 = 3000
    49 = 3000
    56 = 15000
    57 = 2000
    59 = 5000
    64 = 2500
    66 = 5000
    69 = 1000
    70 = 600
    72 = 4000
    96 = 2000
}

$highlighted = @(6, 7, 24, 44, 56, 59, 96)

foreach ($row in $amounts.Keys) {
    $cell = $ws.Cells.Item($row, 14)   # column N = 14
    $cell.Value = $amounts[$row]
    if ($highlighted -contains $row) {
        $cell.Interior.Color = 9660795   # matches existing RGB(7B6993) "paid" fill
    }
}

# Reflect where the user's cursor/scroll ended up after this data entry pass.
$ws.Range("J45").Select()
